$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert five new paragraphs right after the paragraph that ends with
#    " to receiver)" and right before the existing (empty) paragraph that
#    precedes "4. If Sender receive the request...".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("to receiver)")
$rng.Collapse(0)                       # wdCollapseEnd

# New empty paragraph
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)                   # wdCharacter -> step past the new paragraph mark

# "Receiver perspective to sender perspective when resending packets"
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.InsertAfter("Receiver perspective to sender perspective when resending packets")
$rng.Collapse(0)

# "sendTS + latencyWindow +base_Delta > Receiver_now+RTT"
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.InsertAfter("sendTS + latencyWindow +base_Delta > Receiver_now+RTT")
$rng.Collapse(0)

# "sendTS+ latencyWindow+ ½ RTT+clockdiff > Sender_now+clockdiff+ RTT"
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.InsertAfter("sendTS+ latencyWindow+ " + [char]0x00BD + " RTT+clockdiff > Sender_now+clockdiff+ RTT")
$rng.Collapse(0)

# "sendTS+latencyWindow – ½ RTT>Sender_now"
$rng.InsertParagraphAfter()
$rng.Collapse(0)
$rng.MoveStart(1, 1)
$rng.InsertAfter("sendTS+latencyWindow " + [char]0x2013 + " " + [char]0x00BD + " RTT>Sender_now")
$rng.Collapse(0)

# ---------------------------------------------------------------------------
# 2) Mark a rendered page break right before "HalfRTT".
# ---------------------------------------------------------------------------
$lastBreakXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$full = $d.Content.Text
$pos = $full.IndexOf("HalfRTT")
$hRng = $d.Range($pos, $pos)
$hRng.InsertXML($lastBreakXml)

# ---------------------------------------------------------------------------
# 3) Mark a rendered page break right before the final paragraph's leading
#    tab ("      new Base_Delta").
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$pos = $full.LastIndexOf("new Base_Delta")
# back up to the run that begins with the tab character right before the
# run of spaces preceding "new Base_Delta"
$pos = $full.LastIndexOf([char]9, $pos)
$tRng = $d.Range($pos, $pos)
$tRng.InsertXML($lastBreakXml)

Write-Host "Done"
